$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the existing data table and figure out where the new row needs to go.
$table = $ws.ListObjects.Item("Tabela1")
$tableRange = $table.Range
$numCols = $tableRange.Columns.Count
$firstRow = $tableRange.Row
$lastRow = $firstRow + $tableRange.Rows.Count - 1
$newRowNum = $lastRow + 1

# Duplicate the last existing data row so the new row inherits its exact cell
# formatting (date format, number format, alignment, etc.), then overwrite the
# copied values with the new day's data (2020-06-01 / serial 43983).
$ws.Rows.Item($lastRow).Copy()
$ws.Rows.Item($newRowNum).Insert(-4121)

$newRowValues = @(43983, 79698, 659, 1475, 2, 5, 1, 0, 109, 0)
for ($col = 1; $col -le $newRowValues.Length; $col++) {
    $ws.Cells.Item($newRowNum, $col).Value = $newRowValues[$col - 1]
}

# Expand the table (ListObject) so it now includes the newly added row; this
# also grows the table's autoFilter range to match.
$newTableRange = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($newRowNum, $numCols))
$table.Resize($newTableRange)

# Match the selection state left behind after the edit (the newly appended row).
$lastColLetter = [char](64 + $numCols)
$ws.Range("A" + $newRowNum + ":" + $lastColLetter + $newRowNum).Select()
